# Rewrite to cog part 1
# Convert the ALL-CAPS "Types" (column A) and "Cost Group" (column B) values
# in data rows 2..70 into Title Case (first letter of each hyphen-separated
# word capitalized, remainder lower-cased). Header row 1 is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function ConvertTo-TitleCaseHyphen($s) {
    if ([string]::IsNullOrEmpty($s)) { return $s }
    $parts = $s -split ' - '
    $newParts = @()
    foreach ($p in $parts) {
        if ($p.Length -gt 0) {
            $lower = $p.ToLower()
            $newParts += ($lower.Substring(0,1).ToUpper() + $lower.Substring(1))
        } else {
            $newParts += $p
        }
    }
    return [string]::Join(' - ', $newParts)
}

$lastRow = 70

for ($r = 2; $r -le $lastRow; $r++) {
    $aCell = $ws.Cells.Item($r, 1)
    $aVal = $aCell.Value2
    if ($aVal -ne $null -and $aVal -ne "") {
        $aCell.Value = ConvertTo-TitleCaseHyphen $aVal
    }

    $bCell = $ws.Cells.Item($r, 2)
    $bVal = $bCell.Value2
    if ($bVal -ne $null -and $bVal -ne "") {
        $bCell.Value = ConvertTo-TitleCaseHyphen $bVal
    }
}
